# Fruta / hortaliza, semanal
#
# Inserts a new weekly Primera/Segunda price pair for Brocoli
# (Terminal La Palmera de La Serena) dated 2022-02-18, pushing the
# existing rows 543:578 down to 545:580.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 543:544 (everything from the old row 543
# downward shifts down by two rows).
$ws.Range("A543:A544").EntireRow.Insert()

# New "Primera" row (543)
$ws.Cells.Item(543, 1).Value2 = 8
$ws.Cells.Item(543, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(543, 3).Value2 = "Coquimbo"
$ws.Cells.Item(543, 4).Value2 = 44610
$ws.Cells.Item(543, 5).Value2 = 4
$ws.Cells.Item(543, 6).Value2 = 100112023
$ws.Cells.Item(543, 7).Value2 = "Brócoli"
$ws.Cells.Item(543, 8).Value2 = "Sin especificar"
$ws.Cells.Item(543, 9).Value2 = "Primera"
$ws.Cells.Item(543, 10).Value2 = 2500
$ws.Cells.Item(543, 11).Value2 = 850
$ws.Cells.Item(543, 12).Value2 = 900
$ws.Cells.Item(543, 13).Value2 = 875
$ws.Cells.Item(543, 14).Value2 = "`$/unidad"
$ws.Cells.Item(543, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(543, 16).Value2 = 875
$ws.Cells.Item(543, 17).Value2 = 1
$ws.Cells.Item(543, 18).Value2 = "Hortaliza"

# New "Segunda" row (544)
$ws.Cells.Item(544, 1).Value2 = 8
$ws.Cells.Item(544, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(544, 3).Value2 = "Coquimbo"
$ws.Cells.Item(544, 4).Value2 = 44610
$ws.Cells.Item(544, 5).Value2 = 4
$ws.Cells.Item(544, 6).Value2 = 100112023
$ws.Cells.Item(544, 7).Value2 = "Brócoli"
$ws.Cells.Item(544, 8).Value2 = "Sin especificar"
$ws.Cells.Item(544, 9).Value2 = "Segunda"
$ws.Cells.Item(544, 10).Value2 = 1300
$ws.Cells.Item(544, 11).Value2 = 750
$ws.Cells.Item(544, 12).Value2 = 800
$ws.Cells.Item(544, 13).Value2 = 775
$ws.Cells.Item(544, 14).Value2 = "`$/unidad"
$ws.Cells.Item(544, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(544, 16).Value2 = 775
$ws.Cells.Item(544, 17).Value2 = 1
$ws.Cells.Item(544, 18).Value2 = "Hortaliza"
